# CAD-1153: include the filters and templates to show the 3YC assets and request
# Insert three new columns (commitment, commitment start date, commitment end date)
# right before the "Currency" column (Z) on the Data sheet, shifting the remaining
# columns (Currency..USD MSRP) three places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the width of the column the new ones are inserted after (Y = "Adobe User
# Email") so the first new column can inherit it the way Excel normally would.
$adjacentWidth = $ws.Range("Y1").EntireColumn.ColumnWidth

# Insert 3 new blank columns starting at column Z (26), pushing Currency.. to AC..
$ws.Range("Z1:AB1").EntireColumn.Insert() | Out-Null

# Populate the headers for the three newly inserted columns
$ws.Range("Z1").Value = "commitment"
$ws.Range("AA1").Value = "commitment start date"
$ws.Range("AB1").Value = "commitment end date"

# Match Excel's natural behaviour: "commitment" keeps the width of its left
# neighbour, while the two date columns best-fit to their (longer) header text.
$ws.Range("Z1").EntireColumn.ColumnWidth = $adjacentWidth
$ws.Range("AA1:AB1").EntireColumn.AutoFit()

# Refresh the AutoFilter so it spans the new last column (AJ). Toggling it off
# first avoids Range.AutoFilter() just switching an already-active filter off.
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ1").AutoFilter() | Out-Null

# The AutoFilter range change isn't automatically reflected on the workbook's
# hidden "_FilterDatabase" name, so update it explicitly to match.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Data!`$A`$1:`$AJ`$1"
    }
}

# Keep the selection/viewport roughly where Excel would have left it after the edit
$ws.Range("AC9").Select()
